# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'30.029.36"
$ws.Range("E2").Value = "  -0.83%  "

# Row 3
$ws.Range("D3").Value = "'1.902.28"
$ws.Range("E3").Value = "  -1.51%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'0.7403"
$ws.Range("E5").Value = "  -1.60%  "

# Row 6
$ws.Range("D6").Value = "'242.87"
$ws.Range("E6").Value = "  +0.22%  "

# Row 7
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").Value = "'0.3063"
$ws.Range("E8").Value = "  -3.61%  "

# Row 9
$ws.Range("D9").Value = "'25.88"
$ws.Range("E9").Value = "  -6.87%  "

# Row 10
$ws.Range("D10").Value = "'0.06904"
$ws.Range("E10").Value = "  -2.88%  "

# Row 11
$ws.Range("D11").Value = "'0.08011"
$ws.Range("E11").Value = "  -0.44%  "

# Row 12
$ws.Range("D12").Value = "'0.7614"
$ws.Range("E12").Value = "  -2.35%  "

# Row 13
$ws.Range("D13").Value = "'1.902.67"
$ws.Range("E13").Value = "  -1.36%  "

# Row 14
$ws.Range("D14").Value = "'5.242"
$ws.Range("E14").Value = "  -2.79%  "

# Row 15
$ws.Range("D15").Value = "'91.35"
$ws.Range("E15").Value = "  -1.81%  "

# Row 16
$ws.Range("D16").Value = "'6.241"
$ws.Range("E16").Value = "  +3.76%  "

# Row 17
$ws.Range("D17").Value = "'30.033.65"

# Row 18
$ws.Range("D18").Value = "'14.04"
$ws.Range("E18").Value = "  -3.59%  "

# Row 19
$ws.Range("D19").Value = "'0.000007759"
$ws.Range("E19").Value = "  -2.36%  "

# Row 20
$ws.Range("D20").Value = "'237.71"
$ws.Range("E20").Value = "  -5.64%  "

# Row 21
$ws.Range("E21").Value = "  +0.08%  "

# Row 22
$ws.Range("D22").Value = "'2.147.15"
$ws.Range("E22").Value = "  -1.59%  "

# Row 23
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.05%  "

# Row 24
$ws.Range("D24").Value = "'7.053"
$ws.Range("E24").Value = "  +5.42%  "

# Row 25
$ws.Range("D25").Value = "'9.312"
$ws.Range("E25").Value = "  -2.50%  "

# Row 26
$ws.Range("D26").Value = "'166.50"
$ws.Range("E26").Value = "  +0.94%  "

# Row 27
$ws.Range("D27").Value = "'18.82"
$ws.Range("E27").Value = "  -1.47%  "

# Row 28
$ws.Range("D28").Value = "'0.1258"
$ws.Range("E28").Value = "  -3.30%  "

# Row 29
$ws.Range("D29").Value = "'2.041"
$ws.Range("E29").Value = "  -6.76%  "

# Row 30
$ws.Range("E30").Value = "  -1.14%  "

# Row 31
$ws.Range("D31").Value = "'1.534"
$ws.Range("E31").Value = "  -0.75%  "

# Row 32
$ws.Range("E32").Value = "  -2.53%  "

# Row 33
$ws.Range("D33").Value = "'4.050"
$ws.Range("E33").Value = "  -2.33%  "

# Row 34
$ws.Range("D34").Value = "'0.05322"
$ws.Range("E34").Value = "  +1.88%  "

# Row 35
$ws.Range("D35").Value = "'1.295"
$ws.Range("E35").Value = "  -1.65%  "

# Row 36
$ws.Range("D36").Value = "'0.7394"
$ws.Range("E36").Value = "  -2.54%  "

# Row 37
$ws.Range("D37").Value = "'2.728"
$ws.Range("E37").Value = "  -1.93%  "

# Row 38
$ws.Range("D38").Value = "'0.01940"
$ws.Range("E38").Value = "  -0.70%  "

# Row 39
$ws.Range("D39").Value = "'2.794"
$ws.Range("E39").Value = "  -0.12%  "

# Row 40
$ws.Range("E40").Value = "  -3.76%  "

# Row 41
$ws.Range("D41").Value = "'0.4452"
$ws.Range("E41").Value = "  -1.93%  "

# Row 42
$ws.Range("D42").Value = "'73.11"
$ws.Range("E42").Value = "  -6.77%  "

# Row 43
$ws.Range("D43").Value = "'1.964"
$ws.Range("E43").Value = "  -0.81%  "

# Row 44
$ws.Range("E44").Value = "  +0.08%  "

# Row 45
$ws.Range("D45").Value = "'0.8358"
$ws.Range("E45").Value = "  -0.51%  "

# Row 46
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "'7.627"
$ws.Range("E46").Value = "  -0.60%  "

# Row 47
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'101.36"
$ws.Range("E47").Value = "  -0.23%  "

# Row 48
$ws.Range("D48").Value = "'9.790"
$ws.Range("E48").Value = "  -1.75%  "

# Row 49
$ws.Range("D49").Value = "'2.050.44"
$ws.Range("E49").Value = "  -1.66%  "

# Row 50
$ws.Range("E50").Value = "  -3.43%  "

# Row 51
$ws.Range("D51").Value = "'0.1170"
$ws.Range("E51").Value = "  -4.46%  "
